$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.71927011013031
$ws.Range("B1").Value = 1.620960116386414
$ws.Range("C1").Value = 2.040859937667847
$ws.Range("D1").Value = 1.870934963226318
$ws.Range("E1").Value = 2.818394660949707
